# Applies the numeric refresh from the "chore: update Sheets via scheduled runner" commit.
# For each affected row, columns H-N (computed price/profit columns) are updated to their new values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3264.0715
$ws.Range("I51").Value = 2432.8333
$ws.Range("J51").Value = 3887.5
$ws.Range("K51").Value = 2432.8333
$ws.Range("L51").Value = 3887.5
$ws.Range("M51").Value = -1948.8333
$ws.Range("N51").Value = -4855.5

$ws.Range("H53").Value = 1932.25
$ws.Range("I53").Value = 2288.5
$ws.Range("J53").Value = 151
$ws.Range("K53").Value = 2288.5
$ws.Range("L53").Value = 151
$ws.Range("M53").Value = -1651.5
$ws.Range("N53").Value = -1425

$ws.Range("H92").Value = 6206.533
$ws.Range("I92").Value = 3824.75
$ws.Range("J92").Value = 8928.571
$ws.Range("K92").Value = 3824.75
$ws.Range("L92").Value = 8928.571
$ws.Range("M92").Value = -2576.75
$ws.Range("N92").Value = -11424.571

$ws.Range("H137").Value = 11766222
$ws.Range("I137").Value = 1058.3334
$ws.Range("J137").Value = 40002616
$ws.Range("K137").Value = 3175.0002
$ws.Range("L137").Value = 120007848
$ws.Range("M137").Value = -625.0001999999999
$ws.Range("N137").Value = -120012948

$ws.Range("H138").Value = 2004.14
$ws.Range("J138").Value = 2609.3845
$ws.Range("L138").Value = 7828.1535
$ws.Range("N138").Value = -18108.1535

$ws.Range("H139").Value = 39933.332
$ws.Range("J139").Value = 39933.332
$ws.Range("L139").Value = 39933.332
$ws.Range("N139").Value = -50213.332

$ws.Range("H140").Value = 74683.336
$ws.Range("J140").Value = 74683.336
$ws.Range("L140").Value = 74683.336
$ws.Range("N140").Value = -85043.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 2553
$ws.Range("I41").Value = 2553
$ws.Range("K41").Value = 2553
$ws.Range("M41").Value = -2139

$ws.Range("H132").Value = 5815684
$ws.Range("I132").Value = 8929916
$ws.Range("J132").Value = 2451.0667
$ws.Range("K132").Value = 26789748
$ws.Range("L132").Value = 7353.2001
$ws.Range("M132").Value = -26787218
$ws.Range("N132").Value = -12413.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 926.3125
$ws.Range("I99").Value = 831.53845
$ws.Range("J99").Value = 1337
$ws.Range("K99").Value = 831.53845
$ws.Range("L99").Value = 1337
$ws.Range("M99").Value = 666.46155
$ws.Range("N99").Value = -4333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1220.1666
$ws.Range("I16").Value = 1039.1
$ws.Range("J16").Value = 1446.5
$ws.Range("K16").Value = 1039.1
$ws.Range("L16").Value = 1446.5
$ws.Range("M16").Value = -752.0999999999999
$ws.Range("N16").Value = -2020.5

$ws.Range("H31").Value = 8775765
$ws.Range("I31").Value = 4225.212
$ws.Range("J31").Value = 66667932
$ws.Range("K31").Value = 4225.212
$ws.Range("L31").Value = 66667932
$ws.Range("M31").Value = -3930.212
$ws.Range("N31").Value = -66668522

$ws.Range("H34").Value = 8775765
$ws.Range("I34").Value = 4225.212
$ws.Range("J34").Value = 66667932
$ws.Range("K34").Value = 4225.212
$ws.Range("L34").Value = 66667932
$ws.Range("M34").Value = -4023.212
$ws.Range("N34").Value = -66668336

$ws.Range("H58").Value = 3057.838
$ws.Range("I58").Value = 1069.2609
$ws.Range("J58").Value = 6324.7856
$ws.Range("K58").Value = 1069.2609
$ws.Range("L58").Value = 6324.7856
$ws.Range("M58").Value = -866.2609
$ws.Range("N58").Value = -6730.7856

$ws.Range("H99").Value = 1257.4117
$ws.Range("I99").Value = 1146.8334
$ws.Range("K99").Value = 1146.8334
$ws.Range("M99").Value = 351.1666

$ws.Range("H113").Value = 1220.1666
$ws.Range("I113").Value = 1039.1
$ws.Range("J113").Value = 1446.5
$ws.Range("K113").Value = 1039.1
$ws.Range("L113").Value = 1446.5
$ws.Range("M113").Value = 1130.9
$ws.Range("N113").Value = -5786.5

$ws.Range("H122").Value = 1198.0256
$ws.Range("I122").Value = 1276.24
$ws.Range("K122").Value = 3828.72
$ws.Range("M122").Value = -1378.72

$ws.Range("H126").Value = 1257.4117
$ws.Range("I126").Value = 1146.8334
$ws.Range("K126").Value = 3440.5002
$ws.Range("M126").Value = -970.5001999999999

$ws.Range("H132").Value = 2698.8333
$ws.Range("I132").Value = 2020.9524
$ws.Range("J132").Value = 4280.5557
$ws.Range("K132").Value = 6062.857199999999
$ws.Range("L132").Value = 12841.6671
$ws.Range("M132").Value = -3532.857199999999
$ws.Range("N132").Value = -17901.6671

$ws.Range("H136").Value = 3057.838
$ws.Range("I136").Value = 1069.2609
$ws.Range("J136").Value = 6324.7856
$ws.Range("K136").Value = 3207.7827
$ws.Range("L136").Value = 18974.3568
$ws.Range("M136").Value = -657.7826999999997
$ws.Range("N136").Value = -24074.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1530.2778
$ws.Range("I114").Value = 222.66667
$ws.Range("J114").Value = 2837.889
$ws.Range("K114").Value = 668.00001
$ws.Range("L114").Value = 8513.667000000001
$ws.Range("M114").Value = 2585.99999
$ws.Range("N114").Value = -15021.667

$ws.Range("H125").Value = 3414.4443
$ws.Range("J125").Value = 5200
$ws.Range("L125").Value = 15600
$ws.Range("N125").Value = -25440

$ws.Range("H137").Value = 8776019
$ws.Range("I137").Value = 18520082
$ws.Range("J137").Value = 6361.3
$ws.Range("K137").Value = 55560246
$ws.Range("L137").Value = 19083.9
$ws.Range("M137").Value = -55555146
$ws.Range("N137").Value = -29283.9

$ws.Range("H140").Value = 3341.4614
$ws.Range("I140").Value = 2267.182
$ws.Range("J140").Value = 9250
$ws.Range("K140").Value = 6801.545999999999
$ws.Range("L140").Value = 27750
$ws.Range("M140").Value = -1621.545999999999
$ws.Range("N140").Value = -38110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 84583.414
$ws.Range("I113").Value = 101220.1
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 101220.1
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = -99050.10000000001
$ws.Range("N113").Value = -5740

$ws.Range("H126").Value = 3840.6538
$ws.Range("I126").Value = 2364.9167
$ws.Range("K126").Value = 7094.750100000001
$ws.Range("M126").Value = -4624.750100000001

$ws.Range("H138").Value = 57182.668
$ws.Range("J138").Value = 57182.668
$ws.Range("L138").Value = 57182.668
$ws.Range("N138").Value = -67462.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 20098.5
$ws.Range("I74").Value = 20098.5
$ws.Range("K74").Value = 20098.5
$ws.Range("M74").Value = -19100.5

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H77").Value = 20098.5
$ws.Range("I77").Value = 20098.5
$ws.Range("K77").Value = 60295.5
$ws.Range("M77").Value = -55303.5

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H136").Value = 50015276
$ws.Range("I136").Value = 83335460
$ws.Range("J136").Value = 35001.25
$ws.Range("K136").Value = 250006380
$ws.Range("L136").Value = 105003.75
$ws.Range("M136").Value = -250003830
$ws.Range("N136").Value = -110103.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1449.4
$ws.Range("I132").Value = 755.2
$ws.Range("J132").Value = 2143.6
$ws.Range("K132").Value = 2265.6
$ws.Range("L132").Value = 6430.799999999999
$ws.Range("M132").Value = 264.3999999999996
$ws.Range("N132").Value = -11490.8

$ws.Range("H136").Value = 1176.52
$ws.Range("I136").Value = 1175.2
$ws.Range("J136").Value = 1181.8
$ws.Range("K136").Value = 3525.6
$ws.Range("L136").Value = 3545.4
$ws.Range("M136").Value = -975.6000000000004
$ws.Range("N136").Value = -8645.4
